# Update column A on Sheet1: replace the Salesforce metadata *folder*
# names with their corresponding Metadata API *type* names, row for row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @(
    "EmbeddedServiceConfig",
    "EmbeddedServiceLiveAgent",
    "CustomApplication",
    "ApprovalProcess",
    "AuraDefinitionBundle",
    "AutoResponseRules",
    "BrandingSet",
    "CallCenter",
    "ApexClass",
    "Community",
    "ApexComponent",
    "ConnectedApp",
    "ContentAsset",
    "CustomPermission",
    "Dashboard",
    "Document",
    "DuplicateRule",
    "EntitlementProcess",
    "EscalationRules",
    "FlexiPage",
    "Flow",
    "GlobalValueSetTranslation",
    "GlobalValueSet",
    "HomePageLayout",
    "CustomLabels",
    "Layout",
    "Letterhead",
    "LightningExperienceTheme",
    "LightningComponentBundle",
    "MatchingRules",
    "LightningMessageChannel",
    "MilestoneType",
    "Network",
    "NotificationTypeConfig",
    "Translations",
    "CustomObject",
    "ApexPage",
    "PathAssistant",
    "PermissionSet",
    "PresenceUserConfig",
    "QuickAction",
    "RemoteSiteSetting",
    "Role",
    "ServiceChannel",
    "ServicePresenceStatus",
    "SharingCriteriaRule",
    "SiteDotCom",
    "Skill",
    "StandardValueSet",
    "StaticResource",
    "CustomTab",
    "StaticResource",
    "TopicsForObjects",
    "ApexTrigger",
    "Workflow"
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $newValues[$i]
}

# Mirror the author's final cursor position/selection recorded in the
# workbook (scrolled down toward the bottom of the list, cell B58 active).
$ws.Range("B58").Select()
